$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "16/7/2018:" entry.  The paragraph that follows it ends
# with "...controls which gun the player has) etc." and is immediately
# followed by an (already existing) blank paragraph and an "Anis."
# paragraph.  The author continued the journal with a new "18/7/2018:"
# entry right after that "Anis." paragraph, and the "_GoBack" bookmark
# (Word's automatic "last edit location" marker) ends up tracking the
# spot where the new text was typed.
# ------------------------------------------------------------------

$anchor = $d.Content
$anchor.Find.Execute("16/7/2018:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$entryParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($pr.Start -eq $anchor.Start) { $entryParaIndex = $i; break }
}

# Paragraph layout right now (1-based, relative to $entryParaIndex = "16/7/2018:"):
#   entryParaIndex + 0 : "16/7/2018:"
#   entryParaIndex + 1 : "Added a Save functionality ... ) etc."  (contains the _GoBack bookmark)
#   entryParaIndex + 2 : "" (blank)
#   entryParaIndex + 3 : "Anis."
$anisIndex = $entryParaIndex + 3
$anisPara = $d.Paragraphs.Item($anisIndex)

# Remove the old _GoBack bookmark (Word will re-create it at the new
# "last edit" location once we are done typing).
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ------------------------------------------------------------------
# Insert a new blank paragraph after "Anis." (inherits its formatting,
# i.e. non-italic Times New Roman 15pt), then a second new blank
# paragraph after that one (same formatting) which will hold the
# "Cleaned some code..." text.
# ------------------------------------------------------------------
$anisPara.Range.InsertParagraphAfter() | Out-Null
$dateParaIndex = $anisIndex + 1
$datePara = $d.Paragraphs.Item($dateParaIndex)
$datePara.Range.InsertParagraphAfter() | Out-Null

# Fill in the second new paragraph ("Cleaned some code...") FIRST, while
# it is still plain (non-italic), so its rPr/pPr never pick up an <w:i/>
# toggle.
$bodyParaIndex = $dateParaIndex + 1
$bodyPara = $d.Paragraphs.Item($bodyParaIndex)
$bodyText = "Cleaned some code and added NPC Aiming (need to give it a gun, and need to make it aim freely!)"
# Add one throw-away trailing character so the real end-of-text position
# is not (yet) the absolute end of the paragraph -- inserting a zero
# length bookmark exactly at a paragraph boundary snaps it back to the
# start of the document in this host, so we keep one extra character
# past the insertion point until the bookmark is safely placed.
$bodyPara.Range.InsertAfter($bodyText + "X") | Out-Null
$bodyParaAfter = $d.Paragraphs.Item($bodyParaIndex)
$bodyEnd = $bodyParaAfter.Range.End

$bmRange = $d.Range($bodyEnd - 2, $bodyEnd - 2)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$newBookmark = $d.Bookmarks.Item("_GoBack")

$placeholderRange = $d.Range($newBookmark.End, $newBookmark.End + 1)
$placeholderRange.Delete() | Out-Null

# Now fill in the first new paragraph ("18/7/2018:") and italicise it.
$datePara2 = $d.Paragraphs.Item($dateParaIndex)
$datePara2.Range.InsertAfter("18/7/2018:") | Out-Null
$datePara3 = $d.Paragraphs.Item($dateParaIndex)
$datePara3.Range.Font.Italic = 1

Write-Output "done"
